$p = $ppt.ActivePresentation

# 1) Move the "Summary of Team's Contribution" slide (currently #4) to the end
#    of this block, i.e. position #8 (right before the last, "Find My Lecture /
#    Participants" slide). Everything that was at positions 5-8 shifts up by one.
$summarySlide = $p.Slides.Item(4)
$summarySlide.MoveTo(8)

# After the move:
#   position 4 = "Features of the Web Application" (was empty, now gets content)
#   position 5 = "Approach: Front End"            (unchanged, just shifted)
#   position 6 = "Approach: Back End"              (unchanged, just shifted)
#   position 7 = "Live Demonstration"              (unchanged, just shifted)
#   position 8 = "Summary of Team's Contribution"  (was empty, now gets content)

# 2) Fill in the body of "Features of the Web Application" (now slide 4).
$featuresBody = $p.Slides.Item(4).Shapes.Item(2).TextFrame.TextRange
$featuresBody.Text = "Signup – Choose the roll between Professor(Lecturer) and Student"
$featuresBody.InsertAfter("`rLogin, Logout, Reset Password") | Out-Null
$featuresBody.InsertAfter("`rMost viewed lectures and recently uploaded lectures") | Out-Null
$featuresBody.InsertAfter("`rSearch Lecture documents and videos with transcripts") | Out-Null
$featuresBody.InsertAfter("`rFor Students – Saved lectures") | Out-Null
$featuresBody.InsertAfter("`rFor Professors – Uploaded lectures") | Out-Null
$featuresBody.InsertAfter("`rLecture documents (PDF or PPTX)") | Out-Null
$featuresBody.InsertAfter("`rLecture videos with transcripts") | Out-Null

# 3) Fill in the body of "Summary of Team's Contribution" (now slide 8).
$summaryShape = $p.Slides.Item(8).Shapes.Item(2)
$summaryShape.TextFrame.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>
$summaryBody = $summaryShape.TextFrame.TextRange

$summaryBody.Text = "Ishita"
$summaryBody.InsertAfter(" ") | Out-Null
$summaryBody.InsertAfter("Narsiker") | Out-Null
$summaryBody.InsertAfter(" : Backend related with lecture search") | Out-Null

$summaryBody.InsertAfter("`rAmy Eden : Frontend") | Out-Null

$summaryBody.InsertAfter("`rLuke Mullen : Backend related with lectures documents, videos, transcript, and search") | Out-Null

$summaryBody.InsertAfter("`rPetros ") | Out-Null
$summaryBody.InsertAfter("Kitazos") | Out-Null
$summaryBody.InsertAfter(" ") | Out-Null
$summaryBody.InsertAfter(": Frontend") | Out-Null

$summaryBody.InsertAfter("`rSoonKwang Hwang : Backend for reset password, ") | Out-Null
$summaryBody.InsertAfter("darkmode") | Out-Null
